$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 19:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1417398
$ws.Range("C4").Value = 8762
$ws.Range("D4").Value = 300299
$ws.Range("E4").Value = 1033119
$ws.Range("F4").Value = 16406
$ws.Range("G4").Value = 555
$ws.Range("H4").Value = 83980

# Row 9 - Brasil
$ws.Range("B9").Value = 180049
$ws.Range("C9").Value = 2447
$ws.Range("E9").Value = 94853
$ws.Range("G9").Value = 195
$ws.Range("H9").Value = 12599

# Row 15 - India
$ws.Range("B15").Value = 78041
$ws.Range("C15").Value = 3749
$ws.Range("D15").Value = 26300
$ws.Range("E15").Value = 49190
$ws.Range("G15").Value = 136
$ws.Range("H15").Value = 2551

# Row 56 - Marruecos
$ws.Range("B56").Value = 6512
$ws.Range("C56").Value = 94
$ws.Range("D56").Value = 3131
$ws.Range("E56").Value = 3193

# Row 75 - Uzbekistan
$ws.Range("E75").Value = 517
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 11

# Row 92 - Republica de Yibuti
$ws.Range("B92").Value = 1268
$ws.Range("C92").Value = 12
$ws.Range("D92").Value = 900
$ws.Range("E92").Value = 365

# Rows 193/194 - swap Belice / Nueva Caledonia ordering (shared-string swap)
# and their corresponding recovered (D) & death (H) counts
$ws.Range("A193").Value = "Belice"
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2

$ws.Range("A194").Value = "Nueva Caledonia"
$ws.Range("D194").Value = 18
$ws.Range("H194").Value = 0
